$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Nom"
$ws.Range("B1").Value = "Prenom"
$ws.Range("C1").Value = "Matricule"
$ws.Range("D1").Value = "Fonction"
$ws.Range("E1").Value = "Adresse"
$ws.Range("F1").Value = "Date Recrut"
$ws.Range("G1").Value = "Date Detach"
$ws.Range("H1").Value = "Affect Origine"
$ws.Range("I1").Value = "Sit Fam"
$ws.Range("J1").Value = "Nbrs Enfants"

# ---------------------------------------------------------------------------
# Row 2 (BLAL / Mustapha)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "BLAL"
$ws.Range("B2").Value = "Mustapha"
$ws.Range("C2").Value = 224
$ws.Range("D2").Value = "Operateur produit"
$ws.Range("J2").Value = 2

# ---------------------------------------------------------------------------
# Row 3 (ZENDEV / Zoubir)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "ZENDEV"
$ws.Range("B3").Value = "Zoubir"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "Developpeur"
$ws.Range("F3").Value = 45332
$ws.Range("F3").NumberFormat = "m/d/yyyy"
$ws.Range("H3").Value = "Developpeur"
$ws.Range("I3").Value = "veuf"
$ws.Range("E3").Value = "Douera"
$ws.Range("J3").Value = 5

# ---------------------------------------------------------------------------
# Column widths (approximate best-fit sizing for the new columns)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 9.42578125
$ws.Columns.Item(3).ColumnWidth = 9.42578125
$ws.Columns.Item(4).ColumnWidth = 17.42578125
$ws.Columns.Item(5).ColumnWidth = 8.140625
$ws.Columns.Item(6).ColumnWidth = 11.28515625
$ws.Columns.Item(7).ColumnWidth = 11.7109375
$ws.Columns.Item(8).ColumnWidth = 13.5703125
$ws.Columns.Item(10).ColumnWidth = 12.140625

# ---------------------------------------------------------------------------
# Selection (matches the final saved sheet view state)
# ---------------------------------------------------------------------------
$ws.Range("P4").Select()
